$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1252664
$ws.Range("J17").Value = 1252664
$ws.Range("L17").Value = 3757992
$ws.Range("N17").Value = -3758328
$ws.Range("H33").Value = 13159700
$ws.Range("I33").Value = 25000596
$ws.Range("J33").Value = 3149
$ws.Range("K33").Value = 25000596
$ws.Range("L33").Value = 3149
$ws.Range("M33").Value = -25000367
$ws.Range("N33").Value = -3607
$ws.Range("H86").Value = 7639.7
$ws.Range("I86").Value = 10632.333
$ws.Range("J86").Value = 6357.143
$ws.Range("K86").Value = 10632.333
$ws.Range("L86").Value = 6357.143
$ws.Range("M86").Value = -9509.333000000001
$ws.Range("N86").Value = -8603.143
$ws.Range("H89").Value = 7639.7
$ws.Range("I89").Value = 10632.333
$ws.Range("J89").Value = 6357.143
$ws.Range("K89").Value = 53161.665
$ws.Range("L89").Value = 31785.715
$ws.Range("M89").Value = -47545.665
$ws.Range("N89").Value = -43017.715
$ws.Range("H113").Value = 12424.3
$ws.Range("I113").Value = 14820.429
$ws.Range("J113").Value = 6833.3335
$ws.Range("K113").Value = 14820.429
$ws.Range("L113").Value = 6833.3335
$ws.Range("M113").Value = -11566.429
$ws.Range("N113").Value = -13341.3335
$ws.Range("H141").Value = 3168.7273
$ws.Range("I141").Value = 3046.4333
$ws.Range("J141").Value = 4391.6665
$ws.Range("K141").Value = 9139.2999
$ws.Range("L141").Value = 13174.9995
$ws.Range("M141").Value = -3959.2999
$ws.Range("N141").Value = -23534.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 1899.75
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 5000
$ws.Range("N30").Value = -5300
$ws.Range("H45").Value = 1775.4546
$ws.Range("I45").Value = 1481.4445
$ws.Range("K45").Value = 1481.4445
$ws.Range("M45").Value = -1104.4445
$ws.Range("H61").Value = 987669.0600000001
$ws.Range("I61").Value = 3199.25
$ws.Range("K61").Value = 3199.25
$ws.Range("M61").Value = -2987.25
$ws.Range("H74").Value = 60003.855
$ws.Range("J74").Value = 60003.855
$ws.Range("L74").Value = 60003.855
$ws.Range("N74").Value = -61751.855
$ws.Range("H77").Value = 60003.855
$ws.Range("J77").Value = 60003.855
$ws.Range("L77").Value = 300019.275
$ws.Range("N77").Value = -308755.275
$ws.Range("H92").Value = 15500
$ws.Range("J92").Value = 15500
$ws.Range("L92").Value = 15500
$ws.Range("N92").Value = -20492
$ws.Range("H97").Value = 535.36365
$ws.Range("I97").Value = 607.1429000000001
$ws.Range("K97").Value = 607.1429000000001
$ws.Range("M97").Value = -111.1429000000001
$ws.Range("H122").Value = 2580.963
$ws.Range("I122").Value = 1516.375
$ws.Range("K122").Value = 4549.125
$ws.Range("M122").Value = -2099.125
$ws.Range("H136").Value = 987669.0600000001
$ws.Range("I136").Value = 3199.25
$ws.Range("K136").Value = 9597.75
$ws.Range("M136").Value = -7047.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 16193
$ws.Range("I82").Value = 8991.25
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 8991.25
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -8608.25
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 16193
$ws.Range("I85").Value = 8991.25
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 8991.25
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -7665.25
$ws.Range("N85").Value = -47652
$ws.Range("H97").Value = 8200
$ws.Range("I97").Value = 8200
$ws.Range("K97").Value = 8200
$ws.Range("M97").Value = -7209
$ws.Range("H134").Value = 109777.73
$ws.Range("I134").Value = 143653.28
$ws.Range("K134").Value = 430959.84
$ws.Range("M134").Value = -428424.84

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 97352.25999999999
$ws.Range("I31").Value = 183062.19
$ws.Range("J31").Value = 18784.834
$ws.Range("K31").Value = 183062.19
$ws.Range("L31").Value = 18784.834
$ws.Range("M31").Value = -182767.19
$ws.Range("N31").Value = -19374.834
$ws.Range("H34").Value = 97352.25999999999
$ws.Range("I34").Value = 183062.19
$ws.Range("J34").Value = 18784.834
$ws.Range("K34").Value = 183062.19
$ws.Range("L34").Value = 18784.834
$ws.Range("M34").Value = -182860.19
$ws.Range("N34").Value = -19188.834
$ws.Range("H58").Value = 14897.571
$ws.Range("I58").Value = 5867.8335
$ws.Range("K58").Value = 5867.8335
$ws.Range("M58").Value = -5664.8335
$ws.Range("H80").Value = 20833.334
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -17246
$ws.Range("H83").Value = 20833.334
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -56232
$ws.Range("H136").Value = 14897.571
$ws.Range("I136").Value = 5867.8335
$ws.Range("K136").Value = 17603.5005
$ws.Range("M136").Value = -15053.5005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 77085630
$ws.Range("J4").Value = 502000
$ws.Range("L4").Value = 1506000
$ws.Range("N4").Value = -1506224
$ws.Range("H7").Value = 1138.8422
$ws.Range("I7").Value = 201.83333
$ws.Range("J7").Value = 1571.3077
$ws.Range("K7").Value = 605.49999
$ws.Range("L7").Value = 4713.9231
$ws.Range("M7").Value = -493.49999
$ws.Range("N7").Value = -4937.9231
$ws.Range("H11").Value = 2714.2778
$ws.Range("J11").Value = 3193.1667
$ws.Range("L11").Value = 9579.500100000001
$ws.Range("N11").Value = -9859.500100000001
$ws.Range("H22").Value = 3832.5715
$ws.Range("J22").Value = 3295.8
$ws.Range("L22").Value = 9887.400000000001
$ws.Range("N22").Value = -10225.4
$ws.Range("H27").Value = 3832.5715
$ws.Range("J27").Value = 3295.8
$ws.Range("L27").Value = 9887.400000000001
$ws.Range("N27").Value = -10091.4
$ws.Range("H34").Value = 3364.625
$ws.Range("I34").Value = 1813.1428
$ws.Range("J34").Value = 4571.3335
$ws.Range("K34").Value = 5439.428400000001
$ws.Range("L34").Value = 13714.0005
$ws.Range("M34").Value = -5355.428400000001
$ws.Range("N34").Value = -13882.0005
$ws.Range("H35").Value = 200
$ws.Range("I35").Value = 200
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 600
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = -312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 10250
$ws.Range("J15").Value = 10250
$ws.Range("L15").Value = 10250
$ws.Range("N15").Value = -10826
$ws.Range("H81").Value = 10250
$ws.Range("J81").Value = 10250
$ws.Range("L81").Value = 10250
$ws.Range("N81").Value = -12246
$ws.Range("H84").Value = 10250
$ws.Range("J84").Value = 10250
$ws.Range("L84").Value = 30750
$ws.Range("N84").Value = -40734
$ws.Range("H99").Value = 3786
$ws.Range("I99").Value = 3898.2856
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 3898.2856
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -1652.2856
$ws.Range("N99").Value = -7492
$ws.Range("H126").Value = 6895.077
$ws.Range("I126").Value = 8649.125
$ws.Range("J126").Value = 4088.6
$ws.Range("K126").Value = 25947.375
$ws.Range("L126").Value = 12265.8
$ws.Range("M126").Value = -23477.375
$ws.Range("N126").Value = -17205.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1125
$ws.Range("H16").Value = 1121.6428
$ws.Range("I16").Value = 1059.4546
$ws.Range("J16").Value = 1349.6666
$ws.Range("K16").Value = 1059.4546
$ws.Range("L16").Value = 1349.6666
$ws.Range("M16").Value = -889.4546
$ws.Range("N16").Value = -1689.6666
$ws.Range("H22").Value = 4585.5713
$ws.Range("H27").Value = 4585.5713
$ws.Range("H40").Value = 3525.7273
$ws.Range("I40").Value = 3014.75
$ws.Range("J40").Value = 4888.3335
$ws.Range("K40").Value = 3014.75
$ws.Range("L40").Value = 4888.3335
$ws.Range("M40").Value = -2878.75
$ws.Range("N40").Value = -5160.3335
$ws.Range("H82").Value = 2719.25
$ws.Range("I82").Value = 2772.1667
$ws.Range("K82").Value = 2772.1667
$ws.Range("M82").Value = -2411.1667
$ws.Range("H85").Value = 2719.25
$ws.Range("I85").Value = 2772.1667
$ws.Range("K85").Value = 2772.1667
$ws.Range("M85").Value = -1524.1667
$ws.Range("H93").Value = 4944.615
$ws.Range("I93").Value = 5994.7
$ws.Range("J93").Value = 1444.3334
$ws.Range("K93").Value = 5994.7
$ws.Range("L93").Value = 1444.3334
$ws.Range("M93").Value = -4746.7
$ws.Range("N93").Value = -3940.3334
$ws.Range("H95").Value = 50000
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("H132").Value = 1942177.8
$ws.Range("I132").Value = 3168.2964
$ws.Range("J132").Value = 7759206.5
$ws.Range("K132").Value = 9504.889200000001
$ws.Range("L132").Value = 23277619.5
$ws.Range("M132").Value = -6974.889200000001
$ws.Range("N132").Value = -23282679.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2521.7778
$ws.Range("I2").Value = 2860.1333
$ws.Range("K2").Value = 2860.1333
$ws.Range("M2").Value = -2748.1333
$ws.Range("H49").Value = 29992.5
$ws.Range("I49").Value = 29992.5
$ws.Range("K49").Value = 29992.5
$ws.Range("M49").Value = -29762.5
$ws.Range("H107").Value = 1281.4706
$ws.Range("I107").Value = 1394.9286
$ws.Range("K107").Value = 4184.7858
$ws.Range("M107").Value = -2264.7858
$ws.Range("H126").Value = 5566
$ws.Range("I126").Value = 5849.5
$ws.Range("J126").Value = 4999
$ws.Range("K126").Value = 17548.5
$ws.Range("L126").Value = 14997
$ws.Range("M126").Value = -15078.5
$ws.Range("N126").Value = -19937
